$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text (string) representation,
# matching the inline-string cell type used in the source workbook,
# so values such as trailing zeros and "%" signs are preserved exactly.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "297.04"
$ws.Range("E2").Value = "-4.52%"
$ws.Range("D3").Value = "31.56"
$ws.Range("E3").Value = "-1.24%"
$ws.Range("D4").Value = "5.094"
$ws.Range("E4").Value = "-4.71%"
$ws.Range("D5").Value = "0.07484"
$ws.Range("E5").Value = "-1.55%"
$ws.Range("D6").Value = "7.734"
$ws.Range("E6").Value = "-1.35%"
$ws.Range("E7").Value = "4.88%"
$ws.Range("D8").Value = "3.798"
$ws.Range("E8").Value = "2.30%"
$ws.Range("D9").Value = "0.9325"
$ws.Range("E9").Value = "1.07%"
$ws.Range("D10").Value = "0.1702"
$ws.Range("E10").Value = "-1.24%"
$ws.Range("D11").Value = "0.07154"
$ws.Range("E11").Value = "-5.78%"
$ws.Range("D12").Value = "0.07989"
$ws.Range("E12").Value = "-2.13%"
$ws.Range("D13").Value = "0.03015"
$ws.Range("E13").Value = "0.04%"
$ws.Range("D14").Value = "0.09895"
$ws.Range("E14").Value = "0.19%"
$ws.Range("D15").Value = "0.001500"
$ws.Range("E15").Value = "-1.56%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006340"
$ws.Range("E16").Value = "-3.21%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.454"
$ws.Range("E17").Value = "-1.01%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.227"
$ws.Range("E18").Value = "-0.68%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3285"
$ws.Range("E19").Value = "-0.81%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1327"
$ws.Range("E20").Value = "-0.77%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "4.551"
$ws.Range("E21").Value = "7.72%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "0.04652"
$ws.Range("E22").Value = "2.13%"
$ws.Range("E23").Value = "-4.37%"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").Value = "-0.97%"
$ws.Range("D25").Value = "0.004432"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").Value = "0.03%"
$ws.Range("D27").Value = "0.0001875"
$ws.Range("E27").Value = "7.77%"
$ws.Range("D39").Value = "0.01682"
$ws.Range("E39").Value = "0.27%"
$ws.Range("D40").Value = "0.04467"
$ws.Range("E40").Value = "-2.68%"
$ws.Range("D41").Value = "0.007058"
$ws.Range("E41").Value = "-2.62%"
$ws.Range("D42").Value = "0.1327"
$ws.Range("E42").Value = "-2.89%"
$ws.Range("D43").Value = "0.002060"
$ws.Range("E43").Value = "-8.82%"
$ws.Range("D44").Value = "0.01131"
$ws.Range("E44").Value = "-19.92%"
$ws.Range("D45").Value = "0.00006002"
$ws.Range("E45").Value = "-3.12%"
$ws.Range("D46").Value = "1.930"
$ws.Range("E46").Value = "1.96%"
$ws.Range("E47").Value = "-0.21%"
